$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A: rename the setting keys.
$ws.Range("A2").Value = "BroadcastFilesPath"
$ws.Range("A3").Value = "OOHFilePath"
$ws.Range("A4").Value = "VendorFilePath"
$ws.Range("A5").Value = "PrintFilePath"

# Column B: update the corresponding paths.
$ws.Range("B5").Value = "C:\Users\NTMGRM.RPA1\Documents\Print S9\zip\"
$ws.Range("B2").Value = "C:\Users\NTMGRM.RPA1\Documents\UPProjects\Broadcast\Brodcast to S9 data\"
$ws.Range("B3").Value = "C:\Users\NTMGRM.RPA1\Documents\UPProjects\OOH_S9\Data\manifest\zip\"
$ws.Range("B4").Value = "C:\Users\NTMGRM.RPA1\Documents\UPProjects\Vendor S9\Vendor S9 Data\Vendor Manifest Data\zip\"

# Column C description notes are no longer used for these rows.
$ws.Range("C2:C5").ClearContents()

# Rows 6-8 used to hold SharePointURL / SharePointUser / SharePointPass; clear
# the values but keep the existing cell formatting in place.
$ws.Range("A6:B8").ClearContents()

# The SharePointUser cell (B7) carried a mailto hyperlink - remove it.
foreach ($h in $ws.Hyperlinks) {
    $h.Delete()
}

# Update the saved selection to match the author's final cursor position.
$ws.Range("A28").Select()
